$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Premier League")

# Pre-format the C:G numeric/percentage-looking text columns as Text so that
# values like "93%", "3.30", "7.0" are stored as literal strings (matching the
# source workbook, where every data cell is a shared-string, not a number).
$numericRange = $ws.Range("C2:G21")
$numericRange.NumberFormat = "@"

$ws.Range("D2").Value = "7.0"
$ws.Range("E2").Value = "93%"
$ws.Range("F2").Value = "67%"
$ws.Range("G2").Value = "3.30"
$ws.Range("C3").Value = "1.6"
$ws.Range("D3").Value = "8.0"
$ws.Range("E3").Value = "81%"
$ws.Range("F3").Value = "63%"
$ws.Range("G3").Value = "3.30"
$ws.Range("C4").Value = "1.5"
$ws.Range("D4").Value = "7.6"
$ws.Range("E4").Value = "78%"
$ws.Range("F4").Value = "67%"
$ws.Range("G4").Value = "3.37"
$ws.Range("E5").Value = "85%"
$ws.Range("F5").Value = "74%"
$ws.Range("G5").Value = "3.56"
$ws.Range("C6").Value = "2.4"
$ws.Range("D6").Value = "6.2"
$ws.Range("G6").Value = "3.62"
$ws.Range("C7").Value = "2.4"
$ws.Range("E7").Value = "74%"
$ws.Range("F7").Value = "67%"
$ws.Range("G7").Value = "2.81"
$ws.Range("B8").Value = "West Ham"
$ws.Range("C8").Value = "2.4"
$ws.Range("D8").Value = "4.1"
$ws.Range("F8").Value = "66%"
$ws.Range("G8").Value = "3.33"
$ws.Range("C9").Value = "2.1"
$ws.Range("D9").Value = "4.6"
$ws.Range("E9").Value = "82%"
$ws.Range("F9").Value = "74%"
$ws.Range("G9").Value = "3.78"
$ws.Range("B10").Value = "Brighton"
$ws.Range("D10").Value = "5.6"
$ws.Range("E10").Value = "93%"
$ws.Range("F10").Value = "66%"
$ws.Range("G10").Value = "3.44"
$ws.Range("B11").Value = "Wolves"
$ws.Range("C11").Value = "2.5"
$ws.Range("D11").Value = "4.3"
$ws.Range("E11").Value = "82%"
$ws.Range("F11").Value = "66%"
$ws.Range("G11").Value = "3.07"
$ws.Range("B12").Value = "Chelsea"
$ws.Range("C12").Value = "3.2"
$ws.Range("D12").Value = "5.1"
$ws.Range("F12").Value = "61%"
$ws.Range("G12").Value = "3.35"
$ws.Range("E13").Value = "78%"
$ws.Range("F13").Value = "63%"
$ws.Range("D14").Value = "5.9"
$ws.Range("E14").Value = "92%"
$ws.Range("F14").Value = "64%"
$ws.Range("G14").Value = "3.15"
$ws.Range("B15").Value = "Crystal Palace"
$ws.Range("C15").Value = "1.9"
$ws.Range("D15").Value = "4.7"
$ws.Range("E15").Value = "81%"
$ws.Range("F15").Value = "59%"
$ws.Range("G15").Value = "2.93"
$ws.Range("B16").Value = "Brentford"
$ws.Range("C16").Value = "2.3"
$ws.Range("D16").Value = "4.5"
$ws.Range("E16").Value = "85%"
$ws.Range("F16").Value = "70%"
$ws.Range("G16").Value = "3.30"
$ws.Range("B17").Value = "Everton"
$ws.Range("C17").Value = "2.1"
$ws.Range("D17").Value = "4.6"
$ws.Range("E17").Value = "74%"
$ws.Range("F17").Value = "48%"
$ws.Range("G17").Value = "2.44"
$ws.Range("B18").Value = "Nottingham"
$ws.Range("C18").Value = "2.3"
$ws.Range("D18").Value = "3.7"
$ws.Range("E18").Value = "85%"
$ws.Range("F18").Value = "56%"
$ws.Range("G18").Value = "3.07"
$ws.Range("D19").Value = "5.7"
$ws.Range("E19").Value = "85%"
$ws.Range("F19").Value = "73%"
$ws.Range("G19").Value = "3.50"
$ws.Range("E20").Value = "93%"
$ws.Range("F20").Value = "59%"
$ws.Range("G20").Value = "3.15"
$ws.Range("C21").Value = "2.9"
$ws.Range("E21").Value = "89%"
$ws.Range("F21").Value = "66%"

# Restore the default (unstyled) cell style now that the text is committed,
# so the written cells carry no extra number-format styling.
$numericRange.Style = "Normal"
